$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1369.2075
$ws.Range("I132").Value = 1456.3125
$ws.Range("J132").Value = 533
$ws.Range("K132").Value = 4368.9375
$ws.Range("L132").Value = 1599
$ws.Range("M132").Value = -1838.9375
$ws.Range("N132").Value = -6659

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5525.919
$ws.Range("I132").Value = 2117.5881
$ws.Range("J132").Value = 8423
$ws.Range("K132").Value = 6352.7643
$ws.Range("L132").Value = 25269
$ws.Range("M132").Value = -3822.7643
$ws.Range("N132").Value = -30329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 49900
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 49900
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 49900
$ws.Range("N117").Value = -59078
$ws.Range("H118").Value = 57139.2
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 57139.2
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 57139.2
$ws.Range("N118").Value = -60453.2
$ws.Range("H119").Value = 79800
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 79800
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 79800
$ws.Range("N119").Value = -89476
$ws.Range("H120").Value = 9999
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 9999
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 9999
$ws.Range("N120").Value = -19675
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H123").Value = 77032.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 77032.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 77032.5
$ws.Range("N123").Value = -86832.5
$ws.Range("H124").Value = 79800
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 79800
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 79800
$ws.Range("N124").Value = -89620
$ws.Range("H125").Value = 50000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 50000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840
$ws.Range("H126").Value = 30000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 30000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 30000
$ws.Range("N126").Value = -39880
$ws.Range("H127").Value = 55630
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 55630
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 55630
$ws.Range("N127").Value = -65550
$ws.Range("H128").Value = 4700
$ws.Range("I128").Value = 4700
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 14100
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -11610
$ws.Range("H129").Value = 49999
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49999
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H130").Value = 67000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 67000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 67000
$ws.Range("N130").Value = -77040
$ws.Range("H131").Value = 35000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 35000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 35000
$ws.Range("N131").Value = -45080
$ws.Range("H132").Value = 56215
$ws.Range("I132").Value = 20000
$ws.Range("J132").Value = 65268.75
$ws.Range("K132").Value = 20000
$ws.Range("L132").Value = 65268.75
$ws.Range("M132").Value = -14940
$ws.Range("N132").Value = -75388.75
$ws.Range("H133").Value = 49976
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 49976
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 49976
$ws.Range("N133").Value = -60096
$ws.Range("H134").Value = 2017.8918
$ws.Range("I134").Value = 1849
$ws.Range("J134").Value = 2543.3333
$ws.Range("K134").Value = 5547
$ws.Range("L134").Value = 7629.999899999999
$ws.Range("M134").Value = -3012
$ws.Range("N134").Value = -12699.9999
$ws.Range("H135").Value = 48571.43
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 48571.43
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 48571.43
$ws.Range("N135").Value = -58711.43
$ws.Range("H137").Value = 85898.17999999999
$ws.Range("I137").Value = 65000
$ws.Range("J137").Value = 87988
$ws.Range("K137").Value = 65000
$ws.Range("L137").Value = 87988
$ws.Range("M137").Value = -59900
$ws.Range("N137").Value = -98188
$ws.Range("H138").Value = 73700
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 73700
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 73700
$ws.Range("N138").Value = -83980
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 41137.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 41137.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 41137.5
$ws.Range("N140").Value = -51497.5
$ws.Range("H141").Value = 41945
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 41945
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 41945
$ws.Range("N141").Value = -52305

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4751.2856
$ws.Range("I31").Value = 6559.35
$ws.Range("J31").Value = 3107.5908
$ws.Range("K31").Value = 6559.35
$ws.Range("L31").Value = 3107.5908
$ws.Range("M31").Value = -6264.35
$ws.Range("N31").Value = -3697.5908
$ws.Range("H34").Value = 4751.2856
$ws.Range("I34").Value = 6559.35
$ws.Range("J34").Value = 3107.5908
$ws.Range("K34").Value = 6559.35
$ws.Range("L34").Value = 3107.5908
$ws.Range("M34").Value = -6357.35
$ws.Range("N34").Value = -3511.5908
$ws.Range("H129").Value = 45000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 45000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 45000
$ws.Range("N129").Value = -55000
$ws.Range("H130").Value = 66462.86
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 66462.86
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 66462.86
$ws.Range("N130").Value = -76502.86
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H132").Value = 2550.1316
$ws.Range("I132").Value = 2132.4583
$ws.Range("J132").Value = 3266.1428
$ws.Range("K132").Value = 6397.374899999999
$ws.Range("L132").Value = 9798.428400000001
$ws.Range("M132").Value = -3867.374899999999
$ws.Range("N132").Value = -14858.4284
$ws.Range("H133").Value = 45291.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 45291.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 45291.5
$ws.Range("N133").Value = -50351.5
$ws.Range("H134").Value = 3681.8918
$ws.Range("I134").Value = 2328
$ws.Range("J134").Value = 4964.5264
$ws.Range("K134").Value = 6984
$ws.Range("L134").Value = 14893.5792
$ws.Range("M134").Value = -4449
$ws.Range("N134").Value = -19963.5792
$ws.Range("H135").Value = 39248.75
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 39248.75
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 39248.75
$ws.Range("N135").Value = -49388.75
$ws.Range("H137").Value = 56640
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 56640
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 56640
$ws.Range("N137").Value = -66840
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("H140").Value = 54300
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 54300
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 54300
$ws.Range("N140").Value = -64660
$ws.Range("H141").Value = 20296
$ws.Range("I141").Value = 20296
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 20296
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -15116

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2735
$ws.Range("I54").Value = 2000
$ws.Range("J54").Value = 2980
$ws.Range("K54").Value = 6000
$ws.Range("L54").Value = 8940
$ws.Range("M54").Value = -5441
$ws.Range("N54").Value = -10058
$ws.Range("H122").Value = 1040.8889
$ws.Range("I122").Value = 595
$ws.Range("J122").Value = 1118.4348
$ws.Range("K122").Value = 5355
$ws.Range("L122").Value = 10065.9132
$ws.Range("M122").Value = -2905
$ws.Range("N122").Value = -14965.9132

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 41635.5
$ws.Range("J68").Value = 41635.5
$ws.Range("L68").Value = 41635.5
$ws.Range("N68").Value = -43257.5
$ws.Range("H71").Value = 41635.5
$ws.Range("J71").Value = 41635.5
$ws.Range("L71").Value = 124906.5
$ws.Range("N71").Value = -133018.5

Write-Host "Applied all cell updates"
